$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 88
$ws.Range("H88").Value = 3613.3333
$ws.Range("I88").Value = 2618.8
$ws.Range("J88").Value = 3995.8462
$ws.Range("K88").Value = 2618.8
$ws.Range("L88").Value = 3995.8462
$ws.Range("M88").Value = -2212.8
$ws.Range("N88").Value = -4807.8462

# Row 91
$ws.Range("H91").Value = 3613.3333
$ws.Range("I91").Value = 2618.8
$ws.Range("J91").Value = 3995.8462
$ws.Range("K91").Value = 2618.8
$ws.Range("L91").Value = 3995.8462
$ws.Range("M91").Value = -1214.8
$ws.Range("N91").Value = -6803.8462

# Row 132
$ws.Range("H132").Value = 911.5797
$ws.Range("I132").Value = 880.8677
$ws.Range("K132").Value = 2642.6031
$ws.Range("M132").Value = -112.6031000000003

# Row 135
$ws.Range("H135").Value = 337.5946
$ws.Range("I135").Value = 345.58334
$ws.Range("J135").Value = 50
$ws.Range("K135").Value = 3110.25006
$ws.Range("L135").Value = 450
$ws.Range("M135").Value = -575.2500600000003
$ws.Range("N135").Value = -5520

# Row 137
$ws.Range("H137").Value = 1009.6316
$ws.Range("I137").Value = 719.15625
$ws.Range("K137").Value = 2157.46875
$ws.Range("M137").Value = 392.53125

# Row 138
$ws.Range("H138").Value = 2470.926
$ws.Range("I138").Value = 2470.926
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 7412.778
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -2272.778
$ws.Range("N138").ClearContents()

# Row 141
$ws.Range("H141").Value = 1274665.9
$ws.Range("I141").Value = 1556591.8
$ws.Range("K141").Value = 4669775.4
$ws.Range("M141").Value = -4664595.4


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5399.7896
$ws.Range("I32").Value = 4240.2256
$ws.Range("J32").Value = 10535
$ws.Range("K32").Value = 4240.2256
$ws.Range("L32").Value = 10535
$ws.Range("M32").Value = -3953.2256
$ws.Range("N32").Value = -11109

# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

# Row 45
$ws.Range("H45").Value = 1713.3125
$ws.Range("I45").Value = 1531.5714
$ws.Range("K45").Value = 1531.5714
$ws.Range("M45").Value = -1154.5714

# Row 74
$ws.Range("H74").Value = 1181.9062
$ws.Range("I74").Value = 954.8
$ws.Range("K74").Value = 954.8
$ws.Range("M74").Value = -80.79999999999995

# Row 77
$ws.Range("H77").Value = 1181.9062
$ws.Range("I77").Value = 954.8
$ws.Range("K77").Value = 4774
$ws.Range("M77").Value = -406

# Row 122
$ws.Range("H122").Value = 2604.8572
$ws.Range("I122").Value = 2222.3333
$ws.Range("K122").Value = 6666.999899999999
$ws.Range("M122").Value = -4216.999899999999

# Row 132
$ws.Range("H132").Value = 1835.3
$ws.Range("I132").Value = 1279.6538
$ws.Range("J132").Value = 2867.2144
$ws.Range("K132").Value = 3838.9614
$ws.Range("L132").Value = 8601.643199999999
$ws.Range("M132").Value = -1308.9614
$ws.Range("N132").Value = -13661.6432


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 74
$ws.Range("H74").Value = 49999.668
$ws.Range("J74").Value = 49999.668
$ws.Range("L74").Value = 49999.668
$ws.Range("N74").Value = -51871.668

# Row 77
$ws.Range("H77").Value = 49999.668
$ws.Range("J77").Value = 49999.668
$ws.Range("L77").Value = 149999.004
$ws.Range("N77").Value = -159359.004

# Row 80
$ws.Range("H80").Value = 14864.143
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 17258.166
$ws.Range("K80").Value = 500
$ws.Range("L80").Value = 17258.166
$ws.Range("M80").Value = 498
$ws.Range("N80").Value = -19254.166

# Row 83
$ws.Range("H83").Value = 14864.143
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 17258.166
$ws.Range("K83").Value = 2500
$ws.Range("L83").Value = 86290.83
$ws.Range("M83").Value = 2492
$ws.Range("N83").Value = -96274.83

# Row 86
$ws.Range("H86").Value = 170058.17
$ws.Range("I86").Value = 4599.857
$ws.Range("J86").Value = 401699.8
$ws.Range("K86").Value = 4599.857
$ws.Range("L86").Value = 401699.8
$ws.Range("M86").Value = -3476.857
$ws.Range("N86").Value = -403945.8

# Row 89
$ws.Range("H89").Value = 170058.17
$ws.Range("I89").Value = 4599.857
$ws.Range("J89").Value = 401699.8
$ws.Range("K89").Value = 22999.285
$ws.Range("L89").Value = 2008499
$ws.Range("M89").Value = -17383.285
$ws.Range("N89").Value = -2019731

# Row 94
$ws.Range("H94").Value = 622.9231
$ws.Range("I94").Value = 687.25
$ws.Range("K94").Value = 687.25
$ws.Range("M94").Value = -236.25

# Row 134
$ws.Range("H134").Value = 5073.9287
$ws.Range("I134").Value = 6107
$ws.Range("J134").Value = 2491.25
$ws.Range("K134").Value = 18321
$ws.Range("L134").Value = 7473.75
$ws.Range("M134").Value = -15786
$ws.Range("N134").Value = -12543.75


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1919.96
$ws.Range("I31").Value = 1534.7142
$ws.Range("K31").Value = 1534.7142
$ws.Range("M31").Value = -1239.7142

# Row 34
$ws.Range("H34").Value = 1919.96
$ws.Range("I34").Value = 1534.7142
$ws.Range("K34").Value = 1534.7142
$ws.Range("M34").Value = -1332.7142

# Row 122
$ws.Range("H122").Value = 3868.9285
$ws.Range("I122").Value = 2332
$ws.Range("K122").Value = 6996
$ws.Range("M122").Value = -4546

# Row 132
$ws.Range("H132").Value = 1207.591
$ws.Range("I132").Value = 804.13513
$ws.Range("K132").Value = 2412.40539
$ws.Range("M132").Value = 117.5946100000001

# Row 134
$ws.Range("H134").Value = 1352.3269
$ws.Range("I134").Value = 1106.0975
$ws.Range("K134").Value = 3318.2925
$ws.Range("M134").Value = -783.2925000000005


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1053.9445
$ws.Range("J122").Value = 1104.5
$ws.Range("L122").Value = 9940.5
$ws.Range("N122").Value = -14840.5


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4166
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

# Row 83
$ws.Range("H83").Value = 4166
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# Row 132
$ws.Range("H132").Value = 820527.8
$ws.Range("I132").Value = 1426092.6
$ws.Range("K132").Value = 4278277.800000001
$ws.Range("M132").Value = -4275747.800000001


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2875
$ws.Range("J46").Value = 3750
$ws.Range("L46").Value = 3750
$ws.Range("N46").Value = -4126

# Row 93
$ws.Range("H93").Value = 1276.4
$ws.Range("I93").Value = 853.7143
$ws.Range("J93").Value = 2262.6667
$ws.Range("K93").Value = 853.7143
$ws.Range("L93").Value = 2262.6667
$ws.Range("M93").Value = 394.2857
$ws.Range("N93").Value = -4758.6667

# Row 122
$ws.Range("H122").Value = 7021.875
$ws.Range("I122").Value = 6835
$ws.Range("J122").Value = 7333.3335
$ws.Range("K122").Value = 20505
$ws.Range("L122").Value = 22000.0005
$ws.Range("M122").Value = -18055
$ws.Range("N122").Value = -26900.0005

# Row 132
$ws.Range("H132").Value = 1911.4
$ws.Range("I132").Value = 1658.6177
$ws.Range("K132").Value = 4975.8531
$ws.Range("M132").Value = -2445.8531

# Row 136
$ws.Range("H136").Value = 2691.9092
$ws.Range("I136").Value = 1862.45
$ws.Range("J136").Value = 3968
$ws.Range("K136").Value = 5587.35
$ws.Range("L136").Value = 11904
$ws.Range("M136").Value = -3037.35
$ws.Range("N136").Value = -17004


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1050
$ws.Range("J81").Value = 1050
$ws.Range("L81").Value = 2100
$ws.Range("N81").Value = -4222

# Row 84
$ws.Range("H84").Value = 1050
$ws.Range("J84").Value = 1050
$ws.Range("L84").Value = 10500
$ws.Range("N84").Value = -21108

# Row 96
$ws.Range("H96").Value = 12141.857
$ws.Range("I96").Value = 1997
$ws.Range("J96").Value = 16199.8
$ws.Range("K96").Value = 1997
$ws.Range("L96").Value = 16199.8
$ws.Range("M96").Value = -624
$ws.Range("N96").Value = -18945.8

# Row 100
$ws.Range("H100").Value = 701.4545000000001
$ws.Range("I100").Value = 552
$ws.Range("K100").Value = 1104
$ws.Range("M100").Value = -563

# Row 107
$ws.Range("H107").Value = 1845.2
$ws.Range("I107").Value = 1806.5
$ws.Range("K107").Value = 5419.5
$ws.Range("M107").Value = -3499.5

# Row 132
$ws.Range("H132").Value = 2074.9285
$ws.Range("I132").Value = 1248.8823
$ws.Range("K132").Value = 3746.6469
$ws.Range("M132").Value = -1216.6469

# Row 136
$ws.Range("H136").Value = 13551806
$ws.Range("I136").Value = 16836572
$ws.Range("J136").Value = 2145
$ws.Range("K136").Value = 50509716
$ws.Range("L136").Value = 6435
$ws.Range("M136").Value = -50507166
$ws.Range("N136").Value = -11535

